# Split the last "Mapping Table 76" sheet's two source groups ("Débit/Nombre" and
# "Débit/Unité") out into their own sheets ("Mapping Table 77" and "Mapping Table 78"),
# matching the same row/header layout used by every other Mapping Table sheet.

$wb = $excel.ActiveWorkbook

$table76 = $wb.Worksheets.Item("Mapping Table 76")

# --- Build "Mapping Table 77" (3 rows: header + blank + 1 data row) -------------
# Copy an existing 3-row sheet so header styles (s="1"/s="2"), column widths,
# and blank-row formatting are preserved exactly, then overwrite the data row.
$template3 = $wb.Worksheets.Item("Mapping Table 0")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template3.Copy([System.Reflection.Missing]::Value, $lastSheet)
$table77 = $wb.Worksheets.Item($wb.Worksheets.Count)
$table77.Name = "Mapping Table 77"

$table77.Range("A3").Value = "Elément_posologie/Débit/Nombre"
$table77.Range("C3").Value = "not-related-to"
$table77.Range("D3").ClearContents()

# --- Build "Mapping Table 78" (5 rows: header + blank + 3 data rows) ------------
$template5 = $wb.Worksheets.Item("Mapping Table 1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template5.Copy([System.Reflection.Missing]::Value, $lastSheet)
$table78 = $wb.Worksheets.Item($wb.Worksheets.Count)
$table78.Name = "Mapping Table 78"

$table78.Range("A3").Value = "Elément_posologie/Débit/Unité"
$table78.Range("C3").Value = "related-to"
$table78.Range("D3").Value = "Dosage.doseAndRate.rateQuantity.code"

$table78.Range("A4").Value = "Elément_posologie/Débit/Unité"
$table78.Range("C4").Value = "equivalent"
$table78.Range("D4").Value = "Dosage.doseAndRate.rateQuantity.code"

$table78.Range("A5").Value = "Elément_posologie/Débit/Unité"
$table78.Range("C5").Value = "related-to"
$table78.Range("D5").Value = "Dosage.doseAndRate.rateQuantity.unit"

# --- Trim "Mapping Table 76" down to just the "Débit/Nombre" rows (1-5) ---------
$table76.Rows("6:8").Delete()

# --- Bump the build timestamp on the Metadata sheet, like the source commit -----
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-07-08T16:24:33+00:00"
